$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 20
